$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet + workbook window tweaks ---
$ws.Name = "EditIncomeAndExpense"

# --- Header row (row 1): shift "Expected"->out, add ActualResult/Revise/ExpectedResult ---
# (order matters for shared-string table layout: ActualResult, Revise, ExpectedResult, Y, N)
$ws.Range("A1").Value = "Execute"
$ws.Range("B1").Value = "Test Case"
$ws.Range("C1").Value = "Amount"
$ws.Range("F1").Value = "Result"
$ws.Range("E1").Value = "ActualResult"
$ws.Range("G1").Value = "Revise"
$ws.Range("D1").Value = "ExpectedResult"

# --- Row 2 ---
$ws.Range("A2").Value = "Y"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 150

# --- Row 3 ---
$ws.Range("A3").Value = "N"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = -5
$ws.Range("D3").Value = "NaN"

# --- Row 4 ---
$ws.Range("A4").Value = "N"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "five"
$ws.Range("D4").Value = "NaN"

# --- Row 5 ---
$ws.Range("A5").Value = "N"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "1 0"
$ws.Range("D5").Value = "NaN"

# --- Row 6 (C6 stays an empty red-filled cell, D6 gets the Thai error message) ---
$ws.Range("A6").Value = "N"
$ws.Range("B6").Value = 6
$ws.Range("D6").Value = "กรุณากรอกจำนวน"

# --- Styles: column A centered, column B centered+vcentered (also used on header row) ---
$ws.Range("A2:A6").HorizontalAlignment = -4108
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Range("A1:G1").VerticalAlignment = -4108
$ws.Range("B2:B6").HorizontalAlignment = -4108
$ws.Range("B2:B6").VerticalAlignment = -4108

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 8
$ws.Columns("E").ColumnWidth = 9.77

# --- Selection / view ---
$ws.Range("E2").Select()

Write-Output "EditIncomeAndExpense sheet updated"
